$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "before hand" -> "beforehand" and drop the now-stale spelling
#    proofing marks that wrapped the misspelling.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*before hand*") {
        $rng = $p.Range
        $rng.End = $rng.End - 1          # exclude the paragraph mark
        $rng.Text = ""
        $fragment = '<?xml version="1.0" standalone="yes"?>' +
          '<?mso-application progid="Word.Document"?>' +
          '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
          '<w:document xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
          '<w:p><w:r><w:t xml:space="preserve">It is best to print directions, or maps </w:t></w:r>' +
          '<w:r><w:t>beforehand</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> because generally camping areas have little to no reception. </w:t></w:r></w:p>' +
          '</w:body></w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'
        $rng.InsertXML($fragment)
        break
    }
}

# ------------------------------------------------------------------
# 2. Turn the very last (empty) paragraph into a hyperlink pointing
#    at https://wireframe.cc/eJdRJI, styled with the "Hyperlink"
#    character style.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$linkRange = $lastPara.Range
$url = "https://wireframe.cc/eJdRJI"
$d.Hyperlinks.Add($linkRange, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$textRange = $lastPara.Range
$textRange.End = $textRange.End - 1      # exclude the paragraph mark
$textRange.Style = $d.Styles("Hyperlink")

$hlStyle = $d.Styles("Hyperlink")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.QuickStyle = $false

# ------------------------------------------------------------------
# 3. Style catalogue touch-ups that ride along with the hyperlink
#    insertion: stop hiding "Default Paragraph Font" in the style
#    gallery, and register the "Unresolved Mention" character style.
# ------------------------------------------------------------------
$dpf = $d.Styles("DefaultParagraphFont")
$dpf.UnhideWhenUsed = $true

$um = $d.Styles.Add("UnresolvedMention", 2)
$um.NameLocal = "Unresolved Mention"
$um.BaseStyle = $dpf
$um.Priority = 99
$um.UnhideWhenUsed = $true
$um.Font.Color = 6053472
